{"js": "// Change the Docmosis merge field used for the case-management location\n// from `venue_name` to `external_short_name`, e.g.:\n//   <<caseManagementLocation.venue_name>>  ->  <<caseManagementLocation.external_short_name>>\nconst body = context.document.body;\nconst fieldResults = body.search(\".venue_name\", { matchCase: true, matchWholeWord: false });\nfieldResults.load(\"items/text\");\nawait context.sync();\n\nfor (const r of fieldResults.items) {\n  r.insertText(\".external_short_name\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Helper: re-serialize a paragraph through getOoxml()/insertOoxml() so that\n// transient proofing-check markers (<w:proofErr/>) that split a sentence\n// into extra runs are cleared up and the sentence collapses back into a\n// single run, matching what Word does once the proofer re-validates text\n// that has been touched during an edit.\nasync function normalizeParagraphByText(matchText) {\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n\n  const target = paragraphs.items.find((p) => p.text.indexOf(matchText) !== -1);\n  if (!target) {\n    return false;\n  }\n\n  const ooxml = target.getOoxml();\n  await context.sync();\n\n  const match = /<w:p[ >][\\s\\S]*<\\/w:p>/.exec(ooxml.value);\n  if (!match) {\n    return false;\n  }\n  const paragraphXml = match[0];\n\n  const pkg =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n    paragraphXml +\n    \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\n  target.insertOoxml(pkg, Word.InsertLocation.replace);\n  await context.sync();\n  return true;\n}\n\n// \"... early neutral evaluation, mediation and arbitration. ...\" \u2014 the word\n// \"mediation\" previously carried grammar-check markers splitting the\n// sentence into three runs; collapse it back to one clean run.\nawait normalizeParagraphByText(\"the parties must consider settling this litigation\");\n\n// \"Credit hire\" heading \u2014 \"hire\" previously carried grammar-check markers\n// splitting the heading into two runs; collapse it back to one clean run.\nawait normalizeParagraphByText(\"Credit hire\");\n", "ps1": "# Change the Docmosis merge field used for the case-management location\n# from `venue_name` to `external_short_name`, e.g.:\n#   <<caseManagementLocation.venue_name>>  ->  <<caseManagementLocation.external_short_name>>\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \".venue_name\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \".external_short_name\"\n# wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# Helper: re-assert a paragraph's own OOXML (WordOpenXML) back onto itself.\n# This mirrors what Word does once the proofer re-validates a sentence that\n# was touched during an edit: the transient <w:proofErr/> grammar markers\n# that had split the sentence into extra runs are cleared away and the\n# sentence collapses back down into a single clean run.\nfunction Normalize-ParagraphContaining([string]$matchText) {\n    $searchRange = $d.Content\n    $f = $searchRange.Find\n    $f.ClearFormatting()\n    $f.Text = $matchText\n    $found = $f.Execute()\n    if ($found) {\n        $paraRange = $searchRange.Paragraphs(1).Range\n        $xml = $paraRange.WordOpenXML\n        $paraRange.InsertXML($xml) | Out-Null\n    }\n}\n\n# \"... early neutral evaluation, mediation and arbitration. ...\" \u2014 the word\n# \"mediation\" previously carried grammar-check markers splitting the\n# sentence into three runs; collapse it back to one clean run.\nNormalize-ParagraphContaining \"the parties must consider settling this litigation\"\n\n# \"Credit hire\" heading \u2014 \"hire\" previously carried grammar-check markers\n# splitting the heading into two runs; collapse it back to one clean run.\nNormalize-ParagraphContaining \"Credit hire\"\n"}
